$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold + border + centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I2:I12 = 1, J2:J12 mirrors H2:H12
for ($r = 2; $r -le 12; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
